$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("G:G").Insert()
$ws.Range("G1").Value = "d=6"
$ws.Range("G2").Value = 98.03179212897243
$ws.Range("G3").Value = 98.12855343148618
$ws.Range("G4").Value = 98.00292031246215
$ws.Range("G5").Value = 98.07217701378747
$ws.Range("G6").Value = 98.04924168334375
